$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1440.3572
$ws.Range("I99").Value = 354.8889
$ws.Range("J99").Value = 3394.2
$ws.Range("K99").Value = 1064.6667
$ws.Range("L99").Value = 10182.6
$ws.Range("M99").Value = 433.3333
$ws.Range("N99").Value = -13178.6

$ws.Range("H129").Value = 16668413
$ws.Range("I129").Value = 31251108
$ws.Range("J129").Value = 2475.4285
$ws.Range("K129").Value = 93753324
$ws.Range("L129").Value = 7426.2855
$ws.Range("M129").Value = -93748324
$ws.Range("N129").Value = -17426.2855

$ws.Range("H133").Value = 29713.334
$ws.Range("J133").Value = 29713.334
$ws.Range("L133").Value = 29713.334
$ws.Range("N133").Value = -39833.334

$ws.Range("H136").Value = 30440.77
$ws.Range("J136").Value = 30440.77
$ws.Range("L136").Value = 30440.77
$ws.Range("N136").Value = -40640.77

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3527.91
$ws.Range("I32").Value = 2976.3704
$ws.Range("J32").Value = 5879.2104
$ws.Range("K32").Value = 2976.3704
$ws.Range("L32").Value = 5879.2104
$ws.Range("M32").Value = -2689.3704
$ws.Range("N32").Value = -6453.2104

$ws.Range("H37").Value = 24822.8
$ws.Range("I37").Value = 5000
$ws.Range("K37").Value = 5000
$ws.Range("M37").Value = -4727

$ws.Range("H55").Value = 19853
$ws.Range("J55").Value = 19853
$ws.Range("L55").Value = 19853
$ws.Range("N55").Value = -20483

$ws.Range("H132").Value = 1853.8226
$ws.Range("I132").Value = 1283.2826
$ws.Range("K132").Value = 3849.8478
$ws.Range("M132").Value = -1319.8478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 685100.75
$ws.Range("I86").Value = 1233214.2
$ws.Range("J86").Value = 45634.918
$ws.Range("K86").Value = 1233214.2
$ws.Range("L86").Value = 45634.918
$ws.Range("M86").Value = -1232091.2
$ws.Range("N86").Value = -47880.918

$ws.Range("H89").Value = 685100.75
$ws.Range("I89").Value = 1233214.2
$ws.Range("J89").Value = 45634.918
$ws.Range("K89").Value = 6166071
$ws.Range("L89").Value = 228174.59
$ws.Range("M89").Value = -6160455
$ws.Range("N89").Value = -239406.59

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1976.2262
$ws.Range("I31").Value = 1316.8
$ws.Range("J31").Value = 2945.9707
$ws.Range("K31").Value = 1316.8
$ws.Range("L31").Value = 2945.9707
$ws.Range("M31").Value = -1021.8
$ws.Range("N31").Value = -3535.9707

$ws.Range("H34").Value = 1976.2262
$ws.Range("I34").Value = 1316.8
$ws.Range("J34").Value = 2945.9707
$ws.Range("K34").Value = 1316.8
$ws.Range("L34").Value = 2945.9707
$ws.Range("M34").Value = -1114.8
$ws.Range("N34").Value = -3349.9707

$ws.Range("H58").Value = 8198795.5
$ws.Range("I58").Value = 989.78723
$ws.Range("J58").Value = 35720000
$ws.Range("K58").Value = 989.78723
$ws.Range("L58").Value = 35720000
$ws.Range("M58").Value = -786.78723
$ws.Range("N58").Value = -35720406

$ws.Range("H134").Value = 2153.2273
$ws.Range("I134").Value = 832.7646999999999
$ws.Range("J134").Value = 6642.8
$ws.Range("K134").Value = 2498.2941
$ws.Range("L134").Value = 19928.4
$ws.Range("M134").Value = 36.70589999999993
$ws.Range("N134").Value = -24998.4

$ws.Range("H136").Value = 8198795.5
$ws.Range("I136").Value = 989.78723
$ws.Range("J136").Value = 35720000
$ws.Range("K136").Value = 2969.36169
$ws.Range("L136").Value = 107160000
$ws.Range("M136").Value = -419.3616900000002
$ws.Range("N136").Value = -107165100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 1573.25
$ws.Range("I112").Value = 764.3333
$ws.Range("K112").Value = 2292.9999
$ws.Range("M112").Value = -1184.9999

$ws.Range("H133").Value = 4458.75
$ws.Range("I133").Value = 5310
$ws.Range("J133").Value = 3607.5
$ws.Range("K133").Value = 15930
$ws.Range("L133").Value = 10822.5
$ws.Range("M133").Value = -10870
$ws.Range("N133").Value = -20942.5

$ws.Range("H134").Value = 2310.5
$ws.Range("I134").Value = 493.1111
$ws.Range("J134").Value = 4127.8887
$ws.Range("K134").Value = 1479.3333
$ws.Range("L134").Value = 12383.6661
$ws.Range("M134").Value = 3590.6667
$ws.Range("N134").Value = -22523.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5008000
$ws.Range("J7").Value = 2875429
$ws.Range("L7").Value = 2875429
$ws.Range("N7").Value = -2875653

$ws.Range("H8").Value = 5008000
$ws.Range("J8").Value = 2875429
$ws.Range("L8").Value = 2875429
$ws.Range("N8").Value = -2875707

$ws.Range("H12").Value = 2722855.8
$ws.Range("I12").Value = 3300582.2
$ws.Range("J12").Value = 2021330.6
$ws.Range("K12").Value = 3300582.2
$ws.Range("L12").Value = 2021330.6
$ws.Range("M12").Value = -3300442.2
$ws.Range("N12").Value = -2021610.6

$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064

$ws.Range("H126").Value = 3342.9167
$ws.Range("I126").Value = 2244.2856
$ws.Range("J126").Value = 3795.2942
$ws.Range("K126").Value = 6732.8568
$ws.Range("L126").Value = 11385.8826
$ws.Range("M126").Value = -4262.8568
$ws.Range("N126").Value = -16325.8826

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3314.2856
$ws.Range("I7").Value = 1600
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1600
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1488
$ws.Range("N7").Value = -4224

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H55").Value = 1103.7142
$ws.Range("I55").Value = 256
$ws.Range("J55").Value = 1574.6666
$ws.Range("K55").Value = 256
$ws.Range("L55").Value = 1574.6666
$ws.Range("M55").Value = -83
$ws.Range("N55").Value = -1920.6666

$ws.Range("H122").Value = 3120
$ws.Range("I122").Value = 2631.5789
$ws.Range("J122").Value = 4666.6665
$ws.Range("K122").Value = 7894.736699999999
$ws.Range("L122").Value = 13999.9995
$ws.Range("M122").Value = -5444.736699999999
$ws.Range("N122").Value = -18899.9995

$ws.Range("H126").Value = 3314.2856
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -16940

$ws.Range("H136").Value = 2491.25
$ws.Range("I136").Value = 1992.2106
$ws.Range("J136").Value = 3544.7778
$ws.Range("K136").Value = 5976.6318
$ws.Range("L136").Value = 10634.3334
$ws.Range("M136").Value = -3426.6318
$ws.Range("N136").Value = -15734.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 31202
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 31202
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 31202
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -31538

$ws.Range("H124").Value = 32000
$ws.Range("J124").Value = 32000
$ws.Range("L124").Value = 32000
$ws.Range("N124").Value = -41820
